$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.459.86"
$ws.Range("E2").Value = "  +0.78%  "
$ws.Range("D3").Value = "3.229.98"
$ws.Range("E3").Value = "  +1.47%  "
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'604.43"
$ws.Range("E5").Value = "  +1.53%  "
$ws.Range("D6").Value = "'158.09"
$ws.Range("E6").Value = "  +3.19%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "3.230.14"
$ws.Range("E8").Value = "  +1.54%  "
$ws.Range("E9").Value = "  +1.54%  "
$ws.Range("D10").Value = "'0.162"
$ws.Range("E10").Value = "  +1.96%  "
$ws.Range("E11").Value = "  -6.13%  "
$ws.Range("D12").Value = "'0.506"
$ws.Range("E12").Value = "  -1.70%  "
$ws.Range("D13").Value = "'0.0000276"
$ws.Range("E13").Value = "  +3.38%  "
$ws.Range("E14").Value = "  +0.05%  "
$ws.Range("D15").Value = "3.760.81"
$ws.Range("E15").Value = "  +1.53%  "
$ws.Range("D16").Value = "66.583.55"
$ws.Range("E16").Value = "  +0.85%  "
$ws.Range("E17").Value = "  -0.27%  "
$ws.Range("D18").Value = "3.229.12"
$ws.Range("E18").Value = "  +1.47%  "
$ws.Range("E19").Value = "  +1.09%  "
$ws.Range("D20").Value = "'508.97"
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("D21").Value = "'15.28"
$ws.Range("E21").Value = "  -0.48%  "
$ws.Range("D22").Value = "'0.736"
$ws.Range("E22").Value = "  -0.29%  "
$ws.Range("E23").Value = "  +0.37%  "
$ws.Range("D24").Value = "'14.75"
$ws.Range("E24").Value = "  -1.70%  "
$ws.Range("D25").Value = "'84.80"
$ws.Range("E25").Value = "  -0.16%  "
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("D27").Value = "'2.99"
$ws.Range("E27").Value = "  +0.35%  "
$ws.Range("E28").Value = "  -0.67%  "
$ws.Range("E29").Value = "  +4.89%  "
$ws.Range("D30").Value = "'2.98"
$ws.Range("E30").Value = "  +3.75%  "
$ws.Range("E31").Value = "  +0.79%  "
$ws.Range("D32").Value = "'28.20"
$ws.Range("E32").Value = "  +0.28%  "
$ws.Range("E33").Value = "  +0.08%  "
$ws.Range("E34").Value = "  -3.32%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "'0.0989"
$ws.Range("E35").Value = "  +9.84%  "
$ws.Range("B36").Value = "Filecoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D36").Value = "'6.52"
$ws.Range("E36").Value = "  +0.60%  "
$ws.Range("D37").Value = "'512.20"
$ws.Range("E37").Value = "  +6.36%  "
$ws.Range("D38").Value = "'56.09"
$ws.Range("E38").Value = "  +2.47%  "
$ws.Range("E39").Value = "  +17.97%  "
$ws.Range("E40").Value = "  +0.18%  "
$ws.Range("D41").Value = "'3.05"
$ws.Range("E41").Value = "  +7.34%  "
$ws.Range("E42").Value = "  +6.17%  "
$ws.Range("E43").Value = "  -0.40%  "
$ws.Range("D44").Value = "'0.299"
$ws.Range("E44").Value = "  -0.74%  "
$ws.Range("E45").Value = "  +2.74%  "
$ws.Range("D46").Value = "2.876.13"
$ws.Range("E46").Value = "  -0.67%  "
$ws.Range("E47").Value = "  +0.56%  "
$ws.Range("D48").Value = "'2.40"
$ws.Range("E48").Value = "  +4.13%  "
$ws.Range("E49").Value = "  -0.08%  "
$ws.Range("E50").Value = "  -0.03%  "
$ws.Range("D51").Value = "'2.62"
$ws.Range("E51").Value = "  +0.41%  "
